# Update countries & provincias Spain
#
# This script updates the "Pais" worksheet of the COVID countries workbook:
#  1. Refreshes the "last updated" timestamp in A1.
#  2. Fixes the country names for 4 pairs of rows whose case totals changed
#     enough to swap their ranking order (Emiratos Arabes Unidos / Paises
#     Bajos, Corea del Sur / Bosnia y Herzegovina, and a 4-way reshuffle
#     around Guinea / Guayana Francesa / Albania / Croacia).
#  3. Updates the numeric statistics (Casos totales, Nuevos casos, Casos
#     activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#     rows whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Timestamp update
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 14:37"

# 2. Country name corrections (ranking swaps)
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("A46").Value = "Paises Bajos"

$ws.Range("A77").Value = "Bosnia y Herzegovina"
$ws.Range("A78").Value = "Corea del Sur"

$ws.Range("A94").Value = "Croacia"
$ws.Range("A95").Value = "Guinea"
$ws.Range("A96").Value = "Guayana Francesa"
$ws.Range("A97").Value = "Albania"

# 3. Numeric data updates (row -> Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$rowData = @{
    4   = @(6002092, 1727, 3314664, 2503743, 0, 32, 183685)
    28  = @(117988,  246,  114797,  2996,    0, 1,  195)
    45  = @(68511,   491,  59472,   8661,    0, 0,  378)
    46  = @(68114,   0,    0,       0,       0, 0,  6215)
    62  = @(40338,   374,  37116,   2925,    0, 6,  297)
    77  = @(18920,   311,  12480,   5858,    0, 11, 582)
    78  = @(18706,   441,  14461,   3932,    0, 1,  313)
    80  = @(16627,   90,   14763,   1240,    0, 1,  624)
    87  = @(13045,   71,   6594,    5628,    0, 4,  823)
    88  = @(12274,   440,  1209,    10846,   0, 9,  219)
    94  = @(9192,    304,  6595,    2420,    0, 2,  177)
    95  = @(9167,    0,    8150,    960,     0, 0,  57)
    96  = @(8936,    0,    8461,    419,     0, 0,  56)
    97  = @(8927,    0,    4633,    4031,    0, 0,  263)
    139 = @(2087,    5,    1964,    113,     0, 0,  10)
    179 = @(411,     0,    358,     53,      0, 0,  0)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $values[$i]
    }
}
